$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("FR")
$ws.Range("D3").Value = 0.25
$ws.Range("D4").Value = 2.75
$ws.Range("D5").Value = 7.9
$ws.Range("D6").Value = 14.0
$ws.Range("D7").Value = 21.5
$ws.Range("D8").Value = 25.45

$ws = $wb.Worksheets.Item("FR.rolling")
$ws.Range("C3").Value = 0.1
$ws.Range("D3").Value = 0.3
$ws.Range("E3").Value = 0.1
$ws.Range("F3").Value = 0.25
$ws.Range("G3").Value = 0.65
$ws.Range("H3").Value = 0.25
$ws.Range("I3").Value = 0.4
$ws.Range("J3").Value = 0.85
$ws.Range("K3").Value = 0.4
$ws.Range("C4").Value = 0.55
$ws.Range("D4").Value = 3.1
$ws.Range("E4").Value = 0.55
$ws.Range("F4").Value = 1.2
$ws.Range("G4").Value = 6.9
$ws.Range("H4").Value = 1.15
$ws.Range("I4").Value = 2.8
$ws.Range("J4").Value = 14.45
$ws.Range("K4").Value = 2.8
$ws.Range("C5").Value = 1.2
$ws.Range("D5").Value = 9.15
$ws.Range("E5").Value = 1.15
$ws.Range("F5").Value = 2.2
$ws.Range("G5").Value = 15.2
$ws.Range("H5").Value = 2.15
$ws.Range("I5").Value = 5.2
$ws.Range("J5").Value = 27.5
$ws.Range("K5").Value = 5.15
$ws.Range("C6").Value = 1.85
$ws.Range("D6").Value = 16.05
$ws.Range("E6").Value = 1.8
$ws.Range("F6").Value = 2.8
$ws.Range("G6").Value = 24.2
$ws.Range("H6").Value = 2.7
$ws.Range("I6").Value = 6.35
$ws.Range("J6").Value = 38.65
$ws.Range("K6").Value = 6.25
$ws.Range("C7").Value = 2.55
$ws.Range("D7").Value = 23.85
$ws.Range("E7").Value = 2.5
$ws.Range("F7").Value = 3.45
$ws.Range("G7").Value = 31.05
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 6.95
$ws.Range("J7").Value = 46.6
$ws.Range("K7").Value = 6.85
$ws.Range("C8").Value = 2.8
$ws.Range("D8").Value = 27.95
$ws.Range("E8").Value = 2.65
$ws.Range("F8").Value = 3.75
$ws.Range("G8").Value = 36.3
$ws.Range("H8").Value = 3.6
$ws.Range("I8").Value = 7.05
$ws.Range("J8").Value = 51.1
$ws.Range("K8").Value = 6.95

$ws = $wb.Worksheets.Item("highERC")
$ws.Range("C3").Value = 2.65
$ws.Range("E3").Value = 2.65
$ws.Range("F3").Value = 7.1
$ws.Range("H3").Value = 7.1
$ws.Range("I3").Value = 7.95
$ws.Range("K3").Value = 7.95
$ws.Range("C4").Value = 15.6
$ws.Range("E4").Value = 15.6
$ws.Range("F4").Value = 25.15
$ws.Range("H4").Value = 25.15
$ws.Range("I4").Value = 35.65
$ws.Range("K4").Value = 35.65
$ws.Range("C5").Value = 22.7
$ws.Range("E5").Value = 22.6
$ws.Range("F5").Value = 32.4
$ws.Range("H5").Value = 32.3
$ws.Range("I5").Value = 48.2
$ws.Range("K5").Value = 48.15
$ws.Range("C6").Value = 29.55
$ws.Range("E6").Value = 29.05
$ws.Range("F6").Value = 39.35
$ws.Range("H6").Value = 39.2
$ws.Range("I6").Value = 55.45
$ws.Range("K6").Value = 55.05
$ws.Range("C7").Value = 35.0
$ws.Range("E7").Value = 34.05
$ws.Range("F7").Value = 44.15
$ws.Range("H7").Value = 43.75
$ws.Range("I7").Value = 59.75
$ws.Range("K7").Value = 59.4
$ws.Range("C8").Value = 38.35
$ws.Range("E8").Value = 36.95
$ws.Range("F8").Value = 46.6
$ws.Range("H8").Value = 45.7
$ws.Range("I8").Value = 61.4
$ws.Range("K8").Value = 61.0

$ws = $wb.Worksheets.Item("ERCincrease")
$ws.Range("C3").Value = 2.65
$ws.Range("E3").Value = 2.65
$ws.Range("F3").Value = 7.05
$ws.Range("H3").Value = 7.05
$ws.Range("I3").Value = 7.95
$ws.Range("K3").Value = 7.95
$ws.Range("C4").Value = 20.1
$ws.Range("E4").Value = 19.25
$ws.Range("F4").Value = 27.25
$ws.Range("H4").Value = 27.05
$ws.Range("I4").Value = 38.85
$ws.Range("K4").Value = 38.4
$ws.Range("C5").Value = 32.7
$ws.Range("E5").Value = 28.7
$ws.Range("F5").Value = 38.1
$ws.Range("H5").Value = 35.6
$ws.Range("I5").Value = 53.9
$ws.Range("K5").Value = 51.7
$ws.Range("C6").Value = 45.3
$ws.Range("E6").Value = 38.65
$ws.Range("F6").Value = 50.05
$ws.Range("H6").Value = 45.35
$ws.Range("I6").Value = 62.7
$ws.Range("K6").Value = 59.7
$ws.Range("C7").Value = 54.5
$ws.Range("E7").Value = 45.45
$ws.Range("F7").Value = 58.2
$ws.Range("H7").Value = 51.8
$ws.Range("I7").Value = 67.75
$ws.Range("K7").Value = 64.2
$ws.Range("C8").Value = 59.55
$ws.Range("E8").Value = 49.65
$ws.Range("F8").Value = 62.4
$ws.Range("H8").Value = 55.0
$ws.Range("I8").Value = 70.65
$ws.Range("K8").Value = 66.3

$ws = $wb.Worksheets.Item("SummaryReport")
$ws.Range("C2").Value = 27.95
$ws.Range("D2").Value = 36.3
$ws.Range("E2").Value = 51.1
$ws.Range("H2").Value = 2.65
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 6.95
$ws.Range("M2").Value = 2.8
$ws.Range("N2").Value = 3.75
$ws.Range("O2").Value = 7.05
$ws.Range("H3").Value = 36.95
$ws.Range("I3").Value = 45.7
$ws.Range("J3").Value = 61.0
$ws.Range("M3").Value = 38.35
$ws.Range("N3").Value = 46.6
$ws.Range("O3").Value = 61.4
$ws.Range("H4").Value = 49.65
$ws.Range("I4").Value = 55.0
$ws.Range("J4").Value = 66.3
$ws.Range("M4").Value = 59.55
$ws.Range("N4").Value = 62.4
$ws.Range("O4").Value = 70.65

